$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The document currently ends with a trailing empty paragraph (right before
# the sectPr). We need to append, after that paragraph:
#   1. one new empty paragraph
#   2. a bulleted ("ListParagraph") list item about the h2 database
#   3. a bulleted ("ListParagraph") list item about the jdbc dependency
# ---------------------------------------------------------------------------

# Step 1: add a new empty paragraph right after the existing last paragraph.
$r = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$r.InsertParagraphAfter()

# Step 2: add another empty paragraph - this one will become the first list item.
$r = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$r.InsertParagraphAfter()

# Step 3: fill that (now last) paragraph with the first list item's content.
# We use InsertXML on the paragraph's own Range (which includes its paragraph
# mark) so Word re-creates a clean paragraph mark; this always leaves exactly
# one fresh empty trailing paragraph behind, ready for the next item.
$pItem1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$xmlItem1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>I used h2 database for testing, it is in memory database(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pItem1.Range.InsertXML($xmlItem1)

# Step 4: fill the new trailing (last) paragraph with the second list item's content.
$pItem2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$xmlItem2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Apply </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jdbc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dependency to query the database</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pItem2.Range.InsertXML($xmlItem2)

# Step 5: InsertXML left one more fresh empty paragraph after item 2 - the
# diff shows the jdbc paragraph is the very last paragraph in the body (right
# before sectPr), so remove that trailing paragraph mark, merging it away.
$item2Para = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailingMark = $d.Range($item2Para.Range.End - 1, $item2Para.Range.End)
$trailingMark.Delete()

# Step 6: turn both new text paragraphs into a single bulleted list (same
# numId), using the "ListParagraph" style + Word's default bullet template.
$firstListPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$secondListPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$listRange = $d.Range($firstListPara.Range.Start, $secondListPara.Range.End)
$listRange.Style = "ListParagraph"
$bulletTemplate = $word.ListGalleries.Item(1).ListTemplates.Item(1)
$listRange.ListFormat.ApplyListTemplate($bulletTemplate)
